# Weekly Fruta/Hortaliza refresh for "Naranja" sheet.
# The underlying data block (rows 393:492) is shifted down by two rows
# (393:492 -> 395:494) to make room for two brand-new observations at
# the top of the block (rows 393:394); the two oldest rows that fall
# off the bottom of the original range end up as new rows 493:494.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the existing block (A393:T492) before it gets overwritten,
#    then re-write it two rows lower (A395:T494).
$block = $ws.Range("A393:T492").Value2
$ws.Range("A395:T494").Value2 = $block

# The two rows appended at the very bottom (493:494) did not exist
# before, so they need the same date format Excel applies to the rest
# of column D.
$ws.Range("D493").NumberFormat = $ws.Range("D492").NumberFormat
$ws.Range("D494").NumberFormat = $ws.Range("D492").NumberFormat

# 2) Overwrite the freed-up rows 393:394 with the two new observations.
$row393 = New-Object 'object[,]' 1,20
$row393[0,0]  = 7
$row393[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row393[0,2]  = "Ñuble"
$row393[0,3]  = 44889
$row393[0,4]  = 16
$row393[0,5]  = "Fruta"
$row393[0,6]  = 100102
$row393[0,7]  = "Cítricos"
$row393[0,8]  = 100102005
$row393[0,9]  = "Naranja"
$row393[0,10] = "Navel Late"
$row393[0,11] = "Primera"
$row393[0,12] = 60
$row393[0,13] = 9500
$row393[0,14] = 9500
$row393[0,15] = 9500
$row393[0,16] = "`$/bandeja 15 kilos granel"
$row393[0,17] = "Región de O'Higgins"
$row393[0,18] = 633
$row393[0,19] = 15
$ws.Range("A393:T393").Value2 = $row393

$row394 = New-Object 'object[,]' 1,20
$row394[0,0]  = 7
$row394[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row394[0,2]  = "Ñuble"
$row394[0,3]  = 44889
$row394[0,4]  = 16
$row394[0,5]  = "Fruta"
$row394[0,6]  = 100102
$row394[0,7]  = "Cítricos"
$row394[0,8]  = 100102005
$row394[0,9]  = "Naranja"
$row394[0,10] = "Navel Late"
$row394[0,11] = "Segunda"
$row394[0,12] = 120
$row394[0,13] = 8000
$row394[0,14] = 10000
$row394[0,15] = 9000
$row394[0,16] = "`$/bandeja 15 kilos granel"
$row394[0,17] = "Región de O'Higgins"
$row394[0,18] = 600
$row394[0,19] = 15
$ws.Range("A394:T394").Value2 = $row394

# Keep the date formatting consistent on the two newly written rows too.
$ws.Range("D393").NumberFormat = $ws.Range("D392").NumberFormat
$ws.Range("D394").NumberFormat = $ws.Range("D392").NumberFormat
